# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price updates to the Kujata_Profits workbook
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets), matching the author-supplied diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1835.7142
$ws.Range("I29").Value = 266.66666
$ws.Range("J29").Value = 2263.6365
$ws.Range("K29").Value = 799.9999799999999
$ws.Range("L29").Value = 6790.9095
$ws.Range("M29").Value = -518.9999799999999
$ws.Range("N29").Value = -7352.9095

# Row 38
$ws.Range("H38").Value = 141
$ws.Range("I38").Value = 141
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 423
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -51
$ws.Range("N38").ClearContents()

# Row 41
$ws.Range("H41").Value = 15874104
$ws.Range("I41").Value = 23810386
$ws.Range("J41").Value = 1540.2858
$ws.Range("K41").Value = 23810386
$ws.Range("L41").Value = 1540.2858
$ws.Range("M41").Value = -23809946
$ws.Range("N41").Value = -2420.2858

# Row 43
$ws.Range("H43").Value = 13914639
$ws.Range("I43").Value = 50500.5
$ws.Range("J43").Value = 27778778
$ws.Range("K43").Value = 50500.5
$ws.Range("L43").Value = 27778778
$ws.Range("M43").Value = -50431.5
$ws.Range("N43").Value = -27778916

# Row 58
$ws.Range("H58").Value = 315.375
$ws.Range("I58").Value = 315.375
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 946.125
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -796.125
$ws.Range("N58").ClearContents()

# Row 112
$ws.Range("H112").Value = 2143.9285
$ws.Range("J112").Value = 2143.9285
$ws.Range("L112").Value = 6431.7855
$ws.Range("N112").Value = -8647.7855

# Row 138
$ws.Range("H138").Value = 1228.0471
$ws.Range("I138").Value = 802.7234
$ws.Range("J138").Value = 1754.1052
$ws.Range("K138").Value = 2408.1702
$ws.Range("L138").Value = 5262.3156
$ws.Range("M138").Value = 2731.8298
$ws.Range("N138").Value = -15542.3156

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1403.5
$ws.Range("I45").Value = 1484.2
$ws.Range("K45").Value = 1484.2
$ws.Range("M45").Value = -1107.2

# Row 63
$ws.Range("H63").Value = 2371.4285
$ws.Range("J63").Value = 2100
$ws.Range("L63").Value = 2100
$ws.Range("N63").Value = -3472

# Row 66
$ws.Range("H66").Value = 2371.4285
$ws.Range("J66").Value = 2100
$ws.Range("L66").Value = 10500
$ws.Range("N66").Value = -17364

# Row 112
$ws.Range("H112").Value = 10812.667
$ws.Range("J112").Value = 10812.667
$ws.Range("L112").Value = 10812.667
$ws.Range("N112").Value = -13766.667

# Row 119
$ws.Range("H119").Value = 24500
$ws.Range("J119").Value = 24500
$ws.Range("L119").Value = 24500
$ws.Range("N119").Value = -34176

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3902
$ws.Range("I86").Value = 4291.6313
$ws.Range("J86").Value = 2421.4
$ws.Range("K86").Value = 4291.6313
$ws.Range("L86").Value = 2421.4
$ws.Range("M86").Value = -3168.6313
$ws.Range("N86").Value = -4667.4

# Row 89
$ws.Range("H89").Value = 3902
$ws.Range("I89").Value = 4291.6313
$ws.Range("J89").Value = 2421.4
$ws.Range("K89").Value = 21458.1565
$ws.Range("L89").Value = 12107
$ws.Range("M89").Value = -15842.1565
$ws.Range("N89").Value = -23339

# Row 107
$ws.Range("H107").Value = 1747.1666
$ws.Range("I107").Value = 390.33334
$ws.Range("K107").Value = 390.33334
$ws.Range("M107").Value = 1529.66666

# Row 118
$ws.Range("H118").Value = 10000
$ws.Range("J118").Value = 10000
$ws.Range("L118").Value = 10000
$ws.Range("N118").Value = -13314

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 286.36365
$ws.Range("J4").Value = 286.36365
$ws.Range("L4").Value = 286.36365
$ws.Range("N4").Value = -510.36365

# Row 7
$ws.Range("H7").Value = 237.66667
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 237.66667
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 237.66667
$ws.Range("N7").Value = -463.66667
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2200496.2
$ws.Range("J4").Value = 2460797.8
$ws.Range("L4").Value = 7382393.399999999
$ws.Range("N4").Value = -7382617.399999999

# Row 6
$ws.Range("H6").Value = 955.5714
$ws.Range("I6").Value = 229.66667
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 689.00001
$ws.Range("L6").Value = 4500
$ws.Range("M6").Value = -576.00001
$ws.Range("N6").Value = -4726

# Row 87
$ws.Range("H87").Value = 3300
$ws.Range("J87").Value = 3300
$ws.Range("L87").Value = 9900
$ws.Range("N87").Value = -12396

# Row 90
$ws.Range("H90").Value = 3300
$ws.Range("J90").Value = 3300
$ws.Range("L90").Value = 29700
$ws.Range("N90").Value = -42180

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2075.4062
$ws.Range("I102").Value = 2200.4814
$ws.Range("K102").Value = 2200.4814
$ws.Range("M102").Value = -578.4814000000001

# Row 126
$ws.Range("H126").Value = 1985.2307
$ws.Range("I126").Value = 1790.9
$ws.Range("J126").Value = 2633
$ws.Range("K126").Value = 5372.700000000001
$ws.Range("L126").Value = 7899
$ws.Range("M126").Value = -2902.700000000001
$ws.Range("N126").Value = -12839

# Row 132
$ws.Range("H132").Value = 2509.08
$ws.Range("I132").Value = 2143.7896
$ws.Range("K132").Value = 6431.3688
$ws.Range("M132").Value = -3901.3688

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 191123.97
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 52129.75
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 52129.75
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -52353.75

# Row 68
$ws.Range("H68").Value = 1886.6428
$ws.Range("J68").Value = 2319.4
$ws.Range("L68").Value = 2319.4
$ws.Range("N68").Value = -3817.4

# Row 71
$ws.Range("H71").Value = 1886.6428
$ws.Range("J71").Value = 2319.4
$ws.Range("L71").Value = 11597
$ws.Range("N71").Value = -19085

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 110
$ws.Range("H110").Value = 27661
$ws.Range("J110").Value = 26881.334
$ws.Range("L110").Value = 26881.334
$ws.Range("N110").Value = -35061.334

# Row 133
$ws.Range("H133").Value = 35024.75
$ws.Range("J133").Value = 35024.75
$ws.Range("L133").Value = 35024.75
$ws.Range("N133").Value = -40084.75

# Row 140
$ws.Range("H140").Value = 38067.418
$ws.Range("J140").Value = 38067.418
$ws.Range("L140").Value = 38067.418
$ws.Range("N140").Value = -48427.418

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1606.7059
$ws.Range("I96").Value = 1423.6
$ws.Range("K96").Value = 1423.6
$ws.Range("M96").Value = -50.59999999999991

# Row 100
$ws.Range("H100").Value = 2123.75
$ws.Range("I100").Value = 1498.3334
$ws.Range("K100").Value = 2996.6668
$ws.Range("M100").Value = -2455.6668

# Row 119
$ws.Range("H119").Value = 35349
$ws.Range("J119").Value = 35349
$ws.Range("L119").Value = 35349
$ws.Range("N119").Value = -45025
